# DRomics TODO list update:
#  - "short term" sheet: replace the old "85. suivant" placeholder row with
#    the real point 85, and append points 86 and 87 (each with an owner),
#    plus flag point 81 as "ML prototype in the share".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("short term")

# Point 81 (row 16) owner note: "ML" -> "ML prototype in the share"
$ws.Range("B16").Value = "ML prototype in the share"

# Row 21 used to hold the leftover "85. suivant" placeholder; turn it into
# the real point 85 entry, matching the look of the surrounding rows
# (fill style of rows 17/19) and give it an owner in column B.
$ws.Range("A19").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = "85. add an example in the vignette where the need is obvious to work on log scale for bmdplot, sensitivity plot…."
$ws.Range("B21").Value = "ML"

# New row 22: point 86, styled like rows 16/18/20 (fill style used for
# "done" / plain entries).
$ws.Range("A18").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = "86. Add on the Dromics web page an introduction of the functions to help biological interpretation"

# New row 23: point 87, same style family, two-line tall like the other
# long wrapped entries.
$ws.Range("A18").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = "87. change the figs in the cheat sheet using new options, new terminolgy in coherence to the one adopted in DRomicsIntepreter (biological group, experimental level) and add new functions"
$ws.Range("B23").Value = "ML"
$ws.Rows.Item(23).RowHeight = 29

# Owners for rows 22/23 (added after the text so new shared strings are
# appended in the same order the source workbook uses).
$ws.Range("B22").Value = "A et ML"

# Update the sheet view so the newly added rows are visible / selected,
# mirroring where the author left off editing.
$ws.Activate()
$ws.Range("A29").Select() | Out-Null
